# Commit: "updates with supplemental figures"
# Replaces the PSSM score matrix in B2:K21 (10 columns x 20 amino-acid rows)
# with the re-computed values from the supplemental-figures re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(-18.28303407725985,1.920988931703355,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,2.828326410584665,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,2.199829997359945,2.127082871901639,-18.28303407725985,3.463249375437382,-18.28303407725985,1.505131428201203,-18.28303407725985,2.240214866328905,-18.28303407725985),
    @(-18.28303407725985,1.747009298628755,-18.28303407725985,-18.28303407725985,-18.28303407725985,4.321923797968807,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(2.666129352118821,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,1.87386002027386,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(3.771114782688633,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,1.210893670626249,-18.28303407725985,2.163561609599),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,2.806665312865504,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,1.56878357690169),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,2.446162611846605,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,2.031346698639505,1.88090176869559),
    @(-18.28303407725985,-18.28303407725985,1.207450580785363,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,2.065893942840221),
    @(-18.28303407725985,-18.28303407725985,1.212065575614098,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,2.36625455856569,-18.28303407725985),
    @(-18.28303407725985,2.110902022308828,2.407810509928356,-18.28303407725985,-18.28303407725985,-18.28303407725985,1.471959413281034,1.582773184418969,1.504015569651642,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,-18.28303407725985,1.371351966697736,1.262671384163992,1.675859844585171,-18.28303407725985),
    @(-18.28303407725985,-18.28303407725985,1.581987782029485,-18.28303407725985,-18.28303407725985,-18.28303407725985,1.62262468854984,1.676482585971858,-18.28303407725985,-18.28303407725985),
    @(-18.28303407725985,0.9177136965159141,1.435701128486777,-18.28303407725985,3.165236017007677,-18.28303407725985,2.053918069329215,0.9937287283701355,-18.28303407725985,2.227503145296423),
    @(-18.28303407725985,1.03810350641322,-18.28303407725985,1.959224772765467,-18.28303407725985,-18.28303407725985,2.191773983572894,-18.28303407725985,-18.28303407725985,-18.28303407725985)
)

$startRow = 2
$startCol = 2  # column B

for ($r = 0; $r -lt $newValues.Count; $r++) {
    $rowVals = $newValues[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowVals[$c]
    }
}
